$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.05393618525423
$ws.Range("C2").Value = 0.1015048830267489
$ws.Range("D2").Value = 0.0568217256332062
$ws.Range("F2").Value = 1.293916992484043
$ws.Range("G2").Value = 1.015278400264748
$ws.Range("H2").Value = 0.01892298876830645
$ws.Range("I2").Value = 0.02365894106755384
$ws.Range("J2").Value = 0.7225596189526442
$ws.Range("K2").Value = 0.793807992190203
$ws.Range("L2").Value = 0.07482754071707376
$ws.Range("M2").Value = 0.9651993563561234
$ws.Range("N2").Value = 0.2440231907028334
$ws.Range("B3").Value = 0.9185983582577251
$ws.Range("C3").Value = 0.08980804658331465
$ws.Range("D3").Value = 0.05007293875782892
$ws.Range("F3").Value = 1.241514927111808
$ws.Range("G3").Value = 0.974888162831661
$ws.Range("H3").Value = 0.0231285587688147
$ws.Range("I3").Value = 0.02835289009865072
$ws.Range("J3").Value = 0.7082493219598121
$ws.Range("K3").Value = 0.7757196455829103
$ws.Range("L3").Value = 0.06927203260070769
$ws.Range("M3").Value = 0.8397497706764057
$ws.Range("N3").Value = 0.2158408837584034
$ws.Range("B4").Value = 0.8354244366429668
$ws.Range("C4").Value = 0.08273776513983222
$ws.Range("D4").Value = 0.04596162473915655
$ws.Range("F4").Value = 1.209702149176991
$ws.Range("G4").Value = 0.9503939480848942
$ws.Range("H4").Value = 0.02602413506043821
$ws.Range("I4").Value = 0.03159319267907179
$ws.Range("J4").Value = 0.6996887285860822
$ws.Range("K4").Value = 0.7646941314355189
$ws.Range("L4").Value = 0.06581149405728226
$ws.Range("M4").Value = 0.7629016094537349
$ws.Range("N4").Value = 0.1986179014909197
$ws.Range("B5").Value = 0.8012516318770224
$ws.Range("C5").Value = 0.08007658749368574
$ws.Range("D5").Value = 0.04436772199738925
$ws.Range("F5").Value = 1.195488947140518
$ws.Range("G5").Value = 0.9391734148822337
$ws.Range("H5").Value = 0.02728883219413403
$ws.Range("I5").Value = 0.03310260776732354
$ws.Range("J5").Value = 0.6955798000156364
$ws.Range("K5").Value = 0.7592325220292437
$ws.Range("L5").Value = 0.06431441630954282
$ws.Range("M5").Value = 0.7319034212443682
$ws.Range("N5").Value = 0.1917844036064267
$ws.Range("B6").Value = 0.7952652751153266
$ws.Range("C6").Value = 0.07987609203454582
$ws.Range("D6").Value = 0.04419418396412311
$ws.Range("F6").Value = 1.191503482490013
$ws.Range("G6").Value = 0.9357143154072247
$ws.Range("H6").Value = 0.0275137695708173
$ws.Range("I6").Value = 0.0334900303454102
$ws.Range("J6").Value = 0.6940764340432537
$ws.Range("K6").Value = 0.7571214618029245
$ws.Range("L6").Value = 0.06397490421676011
$ws.Range("M6").Value = 0.7271039647848738
$ws.Range("N6").Value = 0.1908555817808946
$ws.Range("B7").Value = 0.8341119558018875
$ws.Range("C7").Value = 0.08335818228533753
$ws.Range("D7").Value = 0.04618811486098195
$ws.Range("F7").Value = 1.205059572339415
$ws.Range("G7").Value = 0.9458734650211369
$ws.Range("H7").Value = 0.02606799640706237
$ws.Range("I7").Value = 0.03195397675520439
$ws.Range("J7").Value = 0.6973849213008378
$ws.Range("K7").Value = 0.7613308697610961
$ws.Range("L7").Value = 0.06554492429280856
$ws.Range("M7").Value = 0.7634290606179945
$ws.Range("N7").Value = 0.1990863427409693
$ws.Range("B8").Value = 1.006138062960019
$ws.Range("C8").Value = 0.09832286486212638
$ws.Range("D8").Value = 0.05481773605929163
$ws.Range("F8").Value = 1.269897079883179
$ws.Range("G8").Value = 0.9955183015189562
$ws.Range("H8").Value = 0.02033652051213652
$ws.Range("I8").Value = 0.02560812320587669
$ws.Range("J8").Value = 0.7146196840976415
$ws.Range("K8").Value = 0.78321910320998
$ws.Range("L8").Value = 0.07259485001556953
$ws.Range("M8").Value = 0.9231572633401015
$ws.Range("N8").Value = 0.2350350933479035
$ws.Range("B9").Value = 1.344752238612273
$ws.Range("C9").Value = 0.1275276023470013
$ws.Range("D9").Value = 0.07158309855564937
$ws.Range("F9").Value = 1.408169331520199
$ws.Range("G9").Value = 1.103394947408219
$ws.Range("H9").Value = 0.011667354266753
$ws.Range("I9").Value = 0.0156270799852658
$ws.Range("J9").Value = 0.7543310634428479
$ws.Range("K9").Value = 0.8327405860530348
$ws.Range("L9").Value = 0.08650329369526233
$ws.Range("M9").Value = 1.236553413215859
$ws.Range("N9").Value = 0.3052849652198404
$ws.Range("B10").Value = 1.591278783380744
$ws.Range("C10").Value = 0.1516391879683141
$ws.Range("D10").Value = 0.08300762557914254
$ws.Range("F10").Value = 1.486223645858615
$ws.Range("G10").Value = 1.162473433372654
$ws.Range("H10").Value = 0.007469334817000739
$ws.Range("I10").Value = 0.01056972312114901
$ws.Range("J10").Value = 0.7735722714300834
$ws.Range("K10").Value = 0.8542167096212339
$ws.Range("L10").Value = 0.09251044529468899
$ws.Range("M10").Value = 1.471055916229915
$ws.Range("N10").Value = 0.3463510258485485
$ws.Range("B11").Value = 1.683257333026489
$ws.Range("C11").Value = 0.177628946612046
$ws.Range("D11").Value = 0.07799455501788088
$ws.Range("F11").Value = 1.316473641385613
$ws.Range("G11").Value = 1.015115034275198
$ws.Range("H11").Value = 0.02588985627803098
$ws.Range("I11").Value = 0.009986810052237338
$ws.Range("J11").Value = 0.6939570424115971
$ws.Range("K11").Value = 0.7445523843805404
$ws.Range("L11").Value = 0.0683892913201003
$ws.Range("M11").Value = 1.597350729599725
$ws.Range("N11").Value = 0.2724890910701987
$ws.Range("B12").Value = 1.709695416241999
$ws.Range("C12").Value = 0.1946290059296274
$ws.Range("D12").Value = 0.07111723041615647
$ws.Range("F12").Value = 1.171939487740261
$ws.Range("G12").Value = 0.8928280346450777
$ws.Range("H12").Value = 0.06485817383290282
$ws.Range("I12").Value = 0.00984212271545104
$ws.Range("J12").Value = 0.6297768453968757
$ws.Range("K12").Value = 0.6596006541145769
$ws.Range("L12").Value = 0.05716427238443167
$ws.Range("M12").Value = 1.653520238026431
$ws.Range("N12").Value = 0.2089404168455644
$ws.Range("B13").Value = 1.685568888761537
$ws.Range("C13").Value = 0.206528050302154
$ws.Range("D13").Value = 0.06284422657017075
$ws.Range("F13").Value = 1.031594064599481
$ws.Range("G13").Value = 0.7759003544468328
$ws.Range("H13").Value = 0.1213370073994469
$ws.Range("I13").Value = 0.01045246782234699
$ws.Range("J13").Value = 0.5701729669105617
$ws.Range("K13").Value = 0.5835498592176478
$ws.Range("L13").Value = 0.0547007131488817
$ws.Range("M13").Value = 1.659867367624059
$ws.Range("N13").Value = 0.1509003484639635
$ws.Range("B14").Value = 1.645582354692863
$ws.Range("C14").Value = 0.2128811844391834
$ws.Range("D14").Value = 0.05656854052020321
$ws.Range("F14").Value = 0.936086613183825
$ws.Range("G14").Value = 0.6971168360226585
$ws.Range("H14").Value = 0.1713912053091065
$ws.Range("I14").Value = 0.01134472520263596
$ws.Range("J14").Value = 0.5307589351586302
$ws.Range("K14").Value = 0.5349050945100693
$ws.Range("L14").Value = 0.0582214156932892
$ws.Range("M14").Value = 1.642047547519411
$ws.Range("N14").Value = 0.1143777589924326
$ws.Range("B15").Value = 1.624800602748593
$ws.Range("C15").Value = 0.2135120839049023
$ws.Range("D15").Value = 0.05481656670778534
$ws.Range("F15").Value = 0.9115831682351399
$ws.Range("G15").Value = 0.6770153179882925
$ws.Range("H15").Value = 0.184236222827451
$ws.Range("I15").Value = 0.0118890606831652
$ws.Range("J15").Value = 0.5211086433031227
$ws.Range("K15").Value = 0.5233064037466413
$ws.Range("L15").Value = 0.05964748642807471
$ws.Range("M15").Value = 1.627363082528802
$ws.Range("N15").Value = 0.1057783490392836
$ws.Range("B16").Value = 1.523214531755059
$ws.Range("C16").Value = 0.2006087005758275
$ws.Range("D16").Value = 0.05193854425650102
$ws.Range("F16").Value = 0.9058574346690236
$ws.Range("G16").Value = 0.6741967576541725
$ws.Range("H16").Value = 0.1732744308134642
$ws.Range("I16").Value = 0.01402173340615054
$ws.Range("J16").Value = 0.5240312304271981
$ws.Range("K16").Value = 0.5284295318582437
$ws.Range("L16").Value = 0.0575657043308464
$ws.Range("M16").Value = 1.524729969867025
$ws.Range("N16").Value = 0.1019356407132577
$ws.Range("B17").Value = 1.466654647737698
$ws.Range("C17").Value = 0.1873675772192769
$ws.Range("D17").Value = 0.05302398279315668
$ws.Range("F17").Value = 0.9510253931361419
$ws.Range("G17").Value = 0.7128560609417605
$ws.Range("H17").Value = 0.136510719505182
$ws.Range("I17").Value = 0.01524152890712749
$ws.Range("J17").Value = 0.5467840843682836
$ws.Range("K17").Value = 0.5573591501081943
$ws.Range("L17").Value = 0.05318118568153984
$ws.Range("M17").Value = 1.455786810701227
$ws.Range("N17").Value = 0.1173620813137006
$ws.Range("B18").Value = 1.442768000500536
$ws.Range("C18").Value = 0.1720903326281302
$ws.Range("D18").Value = 0.05758065434497439
$ws.Range("F18").Value = 1.05263470523289
$ws.Range("G18").Value = 0.7986281959742314
$ws.Range("H18").Value = 0.08387009751580621
$ws.Range("I18").Value = 0.0153214288811947
$ws.Range("J18").Value = 0.5926271385623778
$ws.Range("K18").Value = 0.6157779487065511
$ws.Range("L18").Value = 0.05117804461989151
$ws.Range("M18").Value = 1.407123021942084
$ws.Range("N18").Value = 0.1551839030472379
$ws.Range("B19").Value = 1.444968469845662
$ws.Range("C19").Value = 0.1586285854847631
$ws.Range("D19").Value = 0.06524695838759698
$ws.Range("F19").Value = 1.193542896605308
$ws.Range("G19").Value = 0.9168500700976665
$ws.Range("H19").Value = 0.03784100633729537
$ws.Range("I19").Value = 0.01503440535790457
$ws.Range("J19").Value = 0.6537987041599393
$ws.Range("K19").Value = 0.6952362159741483
$ws.Range("L19").Value = 0.05841930536226947
$ws.Range("M19").Value = 1.38039569159929
$ws.Range("N19").Value = 0.2158324053258127
$ws.Range("B20").Value = 1.523629951748376
$ws.Range("C20").Value = 0.1474369414924297
$ws.Range("D20").Value = 0.08072900242561332
$ws.Range("F20").Value = 1.45092267647081
$ws.Range("G20").Value = 1.132581905052234
$ws.Range("H20").Value = 0.008487356281472902
$ws.Range("I20").Value = 0.0127128666579841
$ws.Range("J20").Value = 0.7611719153544101
$ws.Range("K20").Value = 0.8379541412125633
$ws.Range("L20").Value = 0.0900222776163595
$ws.Range("M20").Value = 1.412396797886174
$ws.Range("N20").Value = 0.3368241886847017
$ws.Range("B21").Value = 1.716886186464365
$ws.Range("C21").Value = 0.1635939676778406
$ws.Range("D21").Value = 0.09133728246643358
$ws.Range("F21").Value = 1.551405060821082
$ws.Range("G21").Value = 1.212767597370899
$ws.Range("H21").Value = 0.005257986681064364
$ws.Range("I21").Value = 0.009133118240967519
$ws.Range("J21").Value = 0.7936980173728472
$ws.Range("K21").Value = 0.8783440756806087
$ws.Range("L21").Value = 0.1008003807128013
$ws.Range("M21").Value = 1.589443235300394
$ws.Range("N21").Value = 0.3859836030836306
$ws.Range("B22").Value = 1.843517480772761
$ws.Range("C22").Value = 0.1742474481359011
$ws.Range("D22").Value = 0.09736276031229352
$ws.Range("F22").Value = 1.611302969213369
$ws.Range("G22").Value = 1.26107262914276
$ws.Range("H22").Value = 0.003709951730089589
$ws.Range("I22").Value = 0.006916162549660321
$ws.Range("J22").Value = 0.8129233060854375
$ws.Range("K22").Value = 0.9020240251403564
$ws.Range("L22").Value = 0.1061682692161625
$ws.Range("M22").Value = 1.705478900605499
$ws.Range("N22").Value = 0.4116879500711406
$ws.Range("B23").Value = 1.776978185991624
$ws.Range("C23").Value = 0.167729769606197
$ws.Range("D23").Value = 0.09384056429281884
$ws.Range("F23").Value = 1.584451026859668
$ws.Range("G23").Value = 1.240315199731526
$ws.Range("H23").Value = 0.004489074226965317
$ws.Range("I23").Value = 0.007706021787228678
$ws.Range("J23").Value = 0.8052135765513526
$ws.Range("K23").Value = 0.8931612860173388
$ws.Range("L23").Value = 0.1035978757662299
$ws.Range("M23").Value = 1.642366769360478
$ws.Range("N23").Value = 0.397269221271273
$ws.Range("B24").Value = 1.523777900533503
$ws.Range("C24").Value = 0.144778301803612
$ws.Range("D24").Value = 0.08108160650510854
$ws.Range("F24").Value = 1.475995178834623
$ws.Range("G24").Value = 1.154973794229335
$ws.Range("H24").Value = 0.008220347371931891
$ws.Range("I24").Value = 0.01210239522089473
$ws.Range("J24").Value = 0.7726761236732074
$ws.Range("K24").Value = 0.8540879016561504
$ws.Range("L24").Value = 0.09338042385681788
$ws.Range("M24").Value = 1.406226299956955
$ws.Range("N24").Value = 0.3440810054683681
$ws.Range("B25").Value = 1.251532923717519
$ws.Range("C25").Value = 0.120707962073638
$ws.Range("D25").Value = 0.06746724896274969
$ws.Range("F25").Value = 1.362390444326948
$ws.Range("G25").Value = 1.066036893306048
$ws.Range("H25").Value = 0.01374683413546379
$ws.Range("I25").Value = 0.01853556551358082
$ws.Range("J25").Value = 0.7393591560659019
$ws.Range("K25").Value = 0.8133747625359007
$ws.Range("L25").Value = 0.08232660908753076
$ws.Range("M25").Value = 1.153225098243013
$ws.Range("N25").Value = 0.287212172680114
